# edit read file excel data
#
# Before: Sheet, get_data_clo (CLO listing)
# After:  Sheet, get_data_course (course listing), get_data_clo (CLO listing)
#
# The existing "get_data_clo" sheet is renamed to "get_data_course" and
# re-populated with the course-listing table; a brand-new "get_data_clo"
# sheet is inserted right after it, carrying the original CLO-listing
# table. The active tab moves to this new "get_data_clo" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Turn the current "get_data_clo" sheet into "get_data_course" and
#    overwrite its data with the course table.
# ---------------------------------------------------------------------
$courseSheet = $wb.Worksheets.Item("get_data_clo")
$courseSheet.Name = "get_data_course"
$courseSheet.Cells.Clear()

# Keep every value as literal text (matches the source data, which stores
# numbers like "4" / "10" as plain strings, not numeric cells).
$courseSheet.Range("A1:G7").NumberFormat = "@"

$courseData = @(
    @("Code",    "Name",                "Program",     "Credits", "Type",      "Status",    "Actions"),
    @("BL2",     "block1",              "Block chain", "4",       "Mandatory", "Completed", "CLOs"),
    @("BLC01",   "Block chain",         "Block chain", "5",       "Mandatory", "Active",    "CLOs"),
    @("BLC02",   "Block chain cơ bản", "Block chain", "0",       "Mandatory", "Updating",  "CLOs"),
    @("KL01111", "aaaa",                "Block chain", "10",      "Mandatory", "Active",    "CLOs"),
    @("PO001",   "Khoá học CN",        "Block chain", "0",       "Mandatory", "Active",    "CLOs"),
    @("as1",     "abc",                 "Block chain", "2",       "Mandatory", "Active",    "CLOs")
)

for ($r = 0; $r -lt $courseData.Length; $r++) {
    $rowValues = $courseData[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $courseSheet.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# ---------------------------------------------------------------------
# 2) Insert a fresh "get_data_clo" sheet right after "get_data_course"
#    and fill it with the original CLO table.
# ---------------------------------------------------------------------
$cloSheet = $wb.Worksheets.Add($null, $courseSheet)
$cloSheet.Name = "get_data_clo"
$cloSheet.Range("A1:H3").NumberFormat = "@"

$cloData = @(
    @("CLO Code", "Description",                         "Category", "Level",    "Weight",  "Assessment Method", "PLO Mapping", "Actions"),
    @("CL001",    "Mục tiêu khoá học Block chain",      "Skill",    "Remember", "10.00%",  "Quiz",              "0",           "PLO"),
    @("CL002",    "Mục tiêu khoá học 2",                "Skill",    "Analyze",  "10.00%",  "Quiz",              "0",           "PLO")
)

for ($r = 0; $r -lt $cloData.Length; $r++) {
    $rowValues = $cloData[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $cloSheet.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# ---------------------------------------------------------------------
# 3) Make the new "get_data_clo" sheet the active tab (was index 1,
#    now index 2 once the extra sheet is inserted).
# ---------------------------------------------------------------------
$cloSheet.Activate()
